# Applies the admiral cheatsheet edit:
#  - slide 1 ("Google Shape;83;p13", shape Id 83): the hyperlinked run
#    "compute_scale" had accidentally been split into two runs
#    ("c" + "ompute_scale"). Re-set it as a single run "compute_scale".
#  - slide 2 ("Google Shape;116;p14", shape Id 46): the
#    derive_vars_crit_flag() description "Derive criterion flag
#    variables (CRITy, CRITyFL(N))." is updated to "Derive criterion
#    flag variables (CRITy, CRITyF(L/N))."

$p = $ppt.ActivePresentation

function Find-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- Slide 1: fix the split "compute_scale" run -----------------------
$s1 = $p.Slides.Item(1)
$sh1 = Find-ShapeById $s1 83
$tr1 = $sh1.TextFrame.TextRange

for ($i = 1; $i -le $tr1.Paragraphs().Count; $i++) {
    $para = $tr1.Paragraphs($i)
    $ptext = $para.Text.TrimEnd("`r", "`n")
    if ($ptext -eq "compute_scale()") {
        # Re-set the "compute_scale" portion (first 13 characters) as a
        # single contiguous run; this merges the stray "c" / "ompute_scale"
        # split while keeping the existing run formatting/hyperlink.
        $run = $para.Characters(1, 13)
        $run.Text = "compute_scale"
    }
}

# --- Slide 2: update the derive_vars_crit_flag() description ----------
$s2 = $p.Slides.Item(2)
$sh2 = Find-ShapeById $s2 46
$tr2 = $sh2.TextFrame.TextRange

for ($i = 1; $i -le $tr2.Paragraphs().Count; $i++) {
    $para = $tr2.Paragraphs($i)
    $ptext = $para.Text.TrimEnd("`r", "`n")
    if ($ptext -eq "Derive criterion flag variables (CRITy, CRITyFL(N)).") {
        # Replace only the trailing "CRITyFL(N))." portion with
        # "CRITyF(L/N)).", leaving "Derive criterion flag variables
        # (CRITy, " as the first run untouched.
        $tail = $para.Characters(41, 12)
        $tail.Text = "CRITyF(L/N))."
    }
}
